$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 457-458, shifting existing rows 457.. down to 459..
$ws.Range("A457:A458").EntireRow.Insert()

# New row 457: Fukumoto / Primera
$ws.Cells.Item(457, 1).Value2 = 4
$ws.Cells.Item(457, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(457, 3).Value2 = "Los Lagos"
$ws.Cells.Item(457, 4).Value2 = 44753
$ws.Cells.Item(457, 5).Value2 = 10
$ws.Cells.Item(457, 6).Value2 = "Fruta"
$ws.Cells.Item(457, 7).Value2 = 100102
$ws.Cells.Item(457, 8).Value2 = "Cítricos"
$ws.Cells.Item(457, 9).Value2 = 100102005
$ws.Cells.Item(457, 10).Value2 = "Naranja"
$ws.Cells.Item(457, 11).Value2 = "Fukumoto"
$ws.Cells.Item(457, 12).Value2 = "Primera"
$ws.Cells.Item(457, 13).Value2 = 200
$ws.Cells.Item(457, 14).Value2 = 14000
$ws.Cells.Item(457, 15).Value2 = 14000
$ws.Cells.Item(457, 16).Value2 = 14000
$ws.Cells.Item(457, 17).Value2 = "$/caja 15 kilos empedrada"
$ws.Cells.Item(457, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(457, 19).Value2 = 933
$ws.Cells.Item(457, 20).Value2 = 15

# New row 458: Fukumoto / Segunda
$ws.Cells.Item(458, 1).Value2 = 4
$ws.Cells.Item(458, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(458, 3).Value2 = "Los Lagos"
$ws.Cells.Item(458, 4).Value2 = 44753
$ws.Cells.Item(458, 5).Value2 = 10
$ws.Cells.Item(458, 6).Value2 = "Fruta"
$ws.Cells.Item(458, 7).Value2 = 100102
$ws.Cells.Item(458, 8).Value2 = "Cítricos"
$ws.Cells.Item(458, 9).Value2 = 100102005
$ws.Cells.Item(458, 10).Value2 = "Naranja"
$ws.Cells.Item(458, 11).Value2 = "Fukumoto"
$ws.Cells.Item(458, 12).Value2 = "Segunda"
$ws.Cells.Item(458, 13).Value2 = 200
$ws.Cells.Item(458, 14).Value2 = 11000
$ws.Cells.Item(458, 15).Value2 = 11000
$ws.Cells.Item(458, 16).Value2 = 11000
$ws.Cells.Item(458, 17).Value2 = "$/caja 15 kilos empedrada"
$ws.Cells.Item(458, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(458, 19).Value2 = 733
$ws.Cells.Item(458, 20).Value2 = 15

Write-Output "done"
